$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert two new rows right after the current last data row (39),
#    pushing the old rows 44/45 (signature block) down to 46/47.
#    Row 39 (CARMEN ALICIA LLINAS HENRIQUEZ / 2508, bottom-border style)
#    keeps its formatting for now - we fix it up below.
# ------------------------------------------------------------------
$ws.Range("B40:J41").EntireRow.Insert()

# ------------------------------------------------------------------
# 2) Fix up borders/formatting:
#    - Row 41 should end up with the special "bottom border" style that
#      row 39 currently has (it becomes the new last row of the table).
#    - Row 39 and Row 40 should get the regular inner-row style that
#      row 38 has.
# ------------------------------------------------------------------
$ws.Range("B39:J39").Copy()
$ws.Range("B41:J41").PasteSpecial(-4122)

$ws.Range("B38:J38").Copy()
$ws.Range("B39:J39").PasteSpecial(-4122)
$ws.Range("B40:J40").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) Populate the data rows.
#    Row 39 keeps its original data (CARMEN, period 2508).
#    Row 40 (new) = EDWIN JIMENEZ AVILA, period 2509.
#    Row 41 (new) = CARMEN ALICIA LLINAS HENRIQUEZ, period 2509.
# ------------------------------------------------------------------
$ws.Range("B39").Value = "CC"
$ws.Range("C39").Value = "36667861"
$ws.Range("D39").Value = "CARMEN ALICIA LLINAS HENRIQUEZ"
$ws.Range("E39").Value = "2508"
$ws.Range("F39").Value = 56000
$ws.Range("G39").Value = 1400000

$ws.Range("B40").Value = "CC"
$ws.Range("C40").Value = "73194983"
$ws.Range("D40").Value = "EDWIN JIMENEZ AVILA"
$ws.Range("E40").Value = "2509"
$ws.Range("F40").Value = 52600
$ws.Range("G40").Value = 1315000

$ws.Range("B41").Value = "CC"
$ws.Range("C41").Value = "36667861"
$ws.Range("D41").Value = "CARMEN ALICIA LLINAS HENRIQUEZ"
$ws.Range("E41").Value = "2509"
$ws.Range("F41").Value = 56000
$ws.Range("G41").Value = 1400000

# ------------------------------------------------------------------
# 4) Update the summary figures at the top of the sheet:
#    - VALOR MORA total (E11): 1209498 -> 1318098
#    - Cant. Periodos (F13): 13 -> 14
# ------------------------------------------------------------------
$ws.Range("E11").Value = 1318098
$ws.Range("F13").Value = 14
